$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "46"
$ws.Range("C9").Characters(27, 9).Text = "11/10/2025"
$ws.Range("C9").Characters(47, 9).Text = "11/16/2025"

# --- Column width updates (I and J widen to match H) ---
$ws.Columns.Item(9).ColumnWidth = 6.71
$ws.Columns.Item(10).ColumnWidth = 6.71

# --- Data table updates rows 15-30 ---
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 20
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 42.857142857142
$ws.Range("M15").Value = 53.846153846153
$ws.Range("N15").Value = -31.03448275862
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 77.777777777777
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 104
$ws.Range("K16").Value = -1.923076923076
$ws.Range("L16").Value = -20.930232558139
$ws.Range("M16").Value = -55.458515283842
$ws.Range("N16").Value = -88.679245283018
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 250
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 10.526315789473
$ws.Range("I17").Value = 227
$ws.Range("J17").Value = 254
$ws.Range("K17").Value = -10.629921259842
$ws.Range("L17").Value = 4.128440366972
$ws.Range("M17").Value = 167.058823529412
$ws.Range("N17").Value = 4.128440366972
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 30.76923076923
$ws.Range("I18").Value = 187
$ws.Range("J18").Value = 191
$ws.Range("K18").Value = -2.094240837696
$ws.Range("L18").Value = -19.742489270386
$ws.Range("M18").Value = -24.596774193548
$ws.Range("N18").Value = -85.356303837118
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 25
$ws.Range("I19").Value = 520
$ws.Range("J19").Value = 489
$ws.Range("K19").Value = 6.339468302658
$ws.Range("L19").Value = -6.642728904847
$ws.Range("M19").Value = 31.645569620253
$ws.Range("N19").Value = -4.059040590405
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 36
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 320
$ws.Range("J20").Value = 341
$ws.Range("K20").Value = -6.158357771261
$ws.Range("L20").Value = 2.236421725239
$ws.Range("M20").Value = 62.43654822335
$ws.Range("N20").Value = -92.626728110599
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 61.904761904761
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = 32.038834951456
$ws.Range("I21").Value = 1379
$ws.Range("J21").Value = 1403
$ws.Range("K21").Value = -1.710620099786
$ws.Range("L21").Value = -5.870307167235
$ws.Range("M21").Value = 17.762596071733
$ws.Range("N21").Value = -81.143169697798
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = "0"
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -100
$ws.Range("I22").Value = 22
$ws.Range("J22").Value = 34
$ws.Range("K22").Value = -35.294117647058
$ws.Range("L22").Value = -43.589743589743
$ws.Range("M22").Value = 4.761904761904
$ws.Range("N22").Value = "***.*"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 49
$ws.Range("J23").Value = 49
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -31.944444444444
$ws.Range("M23").Value = 48.484848484848
$ws.Range("N23").Value = "***.*"
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 7.692307692307
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = -15.853658536585
$ws.Range("I24").Value = 932
$ws.Range("J24").Value = 1049
$ws.Range("K24").Value = -11.153479504289
$ws.Range("L24").Value = -25.796178343949
$ws.Range("M24").Value = 6.271379703534
$ws.Range("N24").Value = "***.*"
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -34.375
$ws.Range("I25").Value = 311
$ws.Range("J25").Value = 447
$ws.Range("K25").Value = -30.425055928411
$ws.Range("L25").Value = -33.118279569892
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -11.904761904761
$ws.Range("I26").Value = 454
$ws.Range("J26").Value = 466
$ws.Range("K26").Value = -2.575107296137
$ws.Range("L26").Value = 5.581395348837
$ws.Range("M26").Value = 25.761772853185
$ws.Range("N26").Value = "***.*"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = -29.032258064516
$ws.Range("L27").Value = -8.333333333333
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"
$ws.Range("C28").Value = "0"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 52
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 48.571428571428
$ws.Range("L28").Value = 20.930232558139
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = "0"
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -100
$ws.Range("I29").Value = 5
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = -16.666666666666
$ws.Range("L29").Value = 25
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -50
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = "0"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -20
$ws.Range("L30").Value = 33.333333333333
$ws.Range("M30").Value = -20
$ws.Range("N30").Value = -60
